$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Every existing data row (2..408) got its "Förändrad" (changed) date
#    bumped by one day: 45202 -> 45203.
$ws.Range("C2:C408").Value = 45203

# 2) The scrape re-ordered the last few entries and added one brand-new
#    entry ("A 46989-2023"). Net effect on the existing rows 405-408 is
#    just an update of the "Beteckning" (A) and "Area (ha)" (G) values;
#    everything else (B/D/E/H:Q/R) is identical between the old and new
#    rows at each of these positions.
$ws.Cells.Item(405, 1).Value = "A 47103-2023"
$ws.Cells.Item(405, 7).Value = 6

$ws.Cells.Item(406, 1).Value = "A 46989-2023"
$ws.Cells.Item(406, 7).Value = 22.9

$ws.Cells.Item(407, 1).Value = "A 47100-2023"
$ws.Cells.Item(407, 7).Value = 1.3

$ws.Cells.Item(408, 1).Value = "A 47002-2023"
$ws.Cells.Item(408, 7).Value = 1.3

# 3) Append a brand-new row 409 carrying the values that used to belong
#    to the old row 407 ("A 47049-2023").
$ws.Cells.Item(409, 1).Value = "A 47049-2023"

$ws.Cells.Item(409, 2).Value = 45201
$ws.Cells.Item(409, 2).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(409, 3).Value = 45203
$ws.Cells.Item(409, 3).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(409, 4).Value = "VÄSTERBOTTENS LÄN"
$ws.Cells.Item(409, 5).Value = "MALÅ"

$ws.Cells.Item(409, 7).Value = 24.6

for ($c = 8; $c -le 17; $c++) {
    $ws.Cells.Item(409, $c).Value = 0
}

$ws.Cells.Item(409, 18).WrapText = $true
